# Apply the edits described by the diff:
# 1. Update header text in B1 (shared string change)
# 2. Update the selection (active cell) on the sheet
# 3. Update various data values in the range B2:F8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header text change
$ws.Range("B1").Value = "% of cases w/ onset 0-1 days before test"

# 3. Data value updates
$ws.Range("D2").Value = 13.7
$ws.Range("F2").Value = 34.6

$ws.Range("B3").Value = 15.2
$ws.Range("C3").Value = 28.4
$ws.Range("D3").Value = 13.8
$ws.Range("E3").Value = 11
$ws.Range("F3").Value = 31.6

$ws.Range("B4").Value = 14.1
$ws.Range("C4").Value = 31.2
$ws.Range("D4").Value = 17.2
$ws.Range("F4").Value = 28.3

$ws.Range("C5").Value = 34.3
$ws.Range("D5").Value = 16.3
$ws.Range("F5").Value = 24.3

$ws.Range("B6").Value = 21.9
$ws.Range("D6").Value = 18.1

$ws.Range("B7").Value = 22.1
$ws.Range("E7").Value = 3.6
$ws.Range("F7").Value = 3.6

$ws.Range("B8").Value = 27.3
$ws.Range("C8").Value = 49.3
$ws.Range("D8").Value = 20.7

# 2. Selection change (active cell moved to C7)
$ws.Range("C7").Select()
